# Add new Kapiti Golf Course hole data to the "Course" worksheet.
# Mirrors a clipboard paste of 9 data rows (holes 1-9) plus a trailing
# blank row, appended right after the existing data (which ended at row 154).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$courseName = "Kapiti Golf Course"

# Hole, Par, Distance for each of the 9 new rows (155-163).
$holes = @(
    @{ Hole = 1; Par = 4; Distance = 269 },
    @{ Hole = 2; Par = 3; Distance = 173 },
    @{ Hole = 3; Par = 4; Distance = 249 },
    @{ Hole = 4; Par = 3; Distance = 168 },
    @{ Hole = 5; Par = 4; Distance = 363 },
    @{ Hole = 6; Par = 4; Distance = 223 },
    @{ Hole = 7; Par = 4; Distance = 262 },
    @{ Hole = 8; Par = 4; Distance = 281 },
    @{ Hole = 9; Par = 4; Distance = 343 }
)

$courseRating = 34
$slopeRating = 102
$lat = -40.9296324816681
$lon = 175.010768880422

$startRow = 155

for ($i = 0; $i -lt $holes.Count; $i++) {
    $r = $startRow + $i
    $data = $holes[$i]

    $ws.Range("A$r").Value = $courseName
    $ws.Range("B$r").Value = $data.Hole
    $ws.Range("C$r").Value = $data.Par
    $ws.Range("D$r").Value = $data.Distance
    $ws.Range("E$r").Value = $courseRating
    $ws.Range("F$r").Value = $slopeRating
    $ws.Range("G$r").Value = $lat
    $ws.Range("H$r").Value = $lon

    $ws.Range("A$r`:G$r").HorizontalAlignment = -4131

    $hCell = $ws.Range("H$r")
    $hCell.HorizontalAlignment = -4131
    $hCell.Font.Name = "Calibri"
    $hCell.Font.Size = 11
    $hCell.Font.ThemeColor = 1
}

# Trailing blank (but formatted) row, as present in the pasted block.
$blankRow = $startRow + $holes.Count
$ws.Range("A$blankRow`:G$blankRow").HorizontalAlignment = -4131
$hBlank = $ws.Range("H$blankRow")
$hBlank.HorizontalAlignment = -4131
$hBlank.Font.Name = "Calibri"
$hBlank.Font.Size = 11
$hBlank.Font.ThemeColor = 1

# Restore the view roughly where it ended up after the paste.
$ws.Range("K149").Select() | Out-Null
